$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sorted)
$ws.Range("B2").Value = 2.986
$ws.Range("C2").Value = 0.003
$ws.Range("D2").Value = 0.44
$ws.Range("E2").Value = 0.037
$ws.Range("F2").Value = 3.161

# Row 3 (Unsorted)
$ws.Range("B3").Value = 3.983
$ws.Range("C3").Value = 1.29
$ws.Range("D3").Value = 0.581
$ws.Range("E3").Value = 0.038
$ws.Range("F3").Value = 0.023

# Row 4 (Reversed)
$ws.Range("B4").Value = 5.003
$ws.Range("C4").Value = 2.545
$ws.Range("D4").Value = 0.88
$ws.Range("E4").Value = 0.051
$ws.Range("F4").Value = 2.419

# Row 5 (Empty)
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0

# Row 6 (Duplicate)
$ws.Range("B6").Value = 2.386
$ws.Range("C6").Value = 0.004
$ws.Range("D6").Value = 0.26
$ws.Range("E6").Value = 0.039
$ws.Range("F6").Value = 1.217

# Remove the side-note commentary cells
$ws.Range("K2").ClearContents()
$ws.Range("K3").ClearContents()

# Remove the leftover scratch table (labels, values, notes)
$ws.Range("L7").ClearContents()
$ws.Range("M6:Q6").ClearContents()
$ws.Range("M7:Q8").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("M11").ClearContents()

# Update selected cell / view
$ws.Range("G4").Select()
